$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.539.65'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '3.088.13'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.56'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.18'
$ws.Range("E6").Value = '  +6.77%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +2.98%  '
$ws.Range("D9").Value = '3.084.43'
$ws.Range("E9").Value = '  -1.54%  '
$ws.Range("E10").Value = '  -1.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.84'
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.56'
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000241'
$ws.Range("E14").Value = '  -2.05%  '
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("E17").Value = '  -2.46%  '
$ws.Range("D18").Value = '63.534.58'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").Value = '3.089.30'
$ws.Range("E19").Value = '  -1.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '476.35'
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.68'
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.719'
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.41'
$ws.Range("E24").Value = '  +3.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.93'
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.30'
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.02'
$ws.Range("E27").Value = '  +2.80%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.38'
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.69'
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.19'
$ws.Range("E32").Value = '  -1.98%  '
$ws.Range("E33").Value = '  +3.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.29'
$ws.Range("E34").Value = '  -1.90%  '
$ws.Range("D35").Value = '0.0₃0848'
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.38'
$ws.Range("E37").Value = '  +5.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.12'
$ws.Range("E38").Value = '  -0.81%  '
$ws.Range("E39").Value = '  -3.38%  '
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.82'
$ws.Range("E41").Value = '  -1.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '443.99'
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.285'
$ws.Range("E43").Value = '  -2.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0364'
$ws.Range("E44").Value = '  -2.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.10'
$ws.Range("E45").Value = '  +1.16%  '
$ws.Range("D47").Value = '2.802.52'
$ws.Range("E47").Value = '  -3.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.36'
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.52'
$ws.Range("E49").Value = '  +5.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.26'
$ws.Range("E51").Value = '  +1.14%  '
